$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 662.5
$ws.Range("I5").Value = 662.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 662.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -547.5
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1067.091
$ws.Range("I18").Value = 739.6667
$ws.Range("J18").Value = 1460
$ws.Range("K18").Value = 739.6667
$ws.Range("L18").Value = 1460
$ws.Range("M18").Value = -455.6667
$ws.Range("N18").Value = -2028

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 489.35715
$ws.Range("I28").Value = 382.14285
$ws.Range("J28").Value = 596.5714
$ws.Range("K28").Value = 382.14285
$ws.Range("L28").Value = 596.5714
$ws.Range("M28").Value = 102.85715
$ws.Range("N28").Value = -1566.5714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 15807.286
$ws.Range("I40").Value = 18108.5
$ws.Range("K40").Value = 18108.5
$ws.Range("M40").Value = -17933.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1478.4286
$ws.Range("I129").Value = 4074.25
$ws.Range("J129").Value = 867.64703
$ws.Range("K129").Value = 12222.75
$ws.Range("L129").Value = 2602.94109
$ws.Range("M129").Value = -7222.75
$ws.Range("N129").Value = -12602.94109

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7108.7104
$ws.Range("I137").Value = 10489.923
$ws.Range("J137").Value = 5350.48
$ws.Range("K137").Value = 31469.769
$ws.Range("L137").Value = 16051.44
$ws.Range("M137").Value = -28919.769
$ws.Range("N137").Value = -21151.44

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2075.4285
$ws.Range("I138").Value = 1691.1333
$ws.Range("J138").Value = 3036.1667
$ws.Range("K138").Value = 5073.3999
$ws.Range("L138").Value = 9108.500100000001
$ws.Range("M138").Value = 66.60009999999966
$ws.Range("N138").Value = -19388.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 3968.5
$ws.Range("I35").Value = 3968.5
$ws.Range("K35").Value = 3968.5
$ws.Range("M35").Value = -3562.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4374.65
$ws.Range("I132").Value = 3528.5715
$ws.Range("J132").Value = 4830.231
$ws.Range("K132").Value = 10585.7145
$ws.Range("L132").Value = 14490.693
$ws.Range("M132").Value = -8055.7145
$ws.Range("N132").Value = -19550.693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2144.1667
$ws.Range("I134").Value = 1764.6111
$ws.Range("J134").Value = 3282.8333
$ws.Range("K134").Value = 5293.8333
$ws.Range("L134").Value = 9848.499899999999
$ws.Range("M134").Value = -2758.8333
$ws.Range("N134").Value = -14918.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 64936.145
$ws.Range("J137").Value = 64936.145
$ws.Range("L137").Value = 64936.145
$ws.Range("N137").Value = -75136.14499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7003333.5
$ws.Range("I6").Value = 7003333.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 7003333.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -7003220.5
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5054483
$ws.Range("I31").Value = 1858.591
$ws.Range("J31").Value = 7580795
$ws.Range("K31").Value = 1858.591
$ws.Range("L31").Value = 7580795
$ws.Range("M31").Value = -1563.591
$ws.Range("N31").Value = -7581385

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5054483
$ws.Range("I34").Value = 1858.591
$ws.Range("J34").Value = 7580795
$ws.Range("K34").Value = 1858.591
$ws.Range("L34").Value = 7580795
$ws.Range("M34").Value = -1656.591
$ws.Range("N34").Value = -7581199

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1743
$ws.Range("I58").Value = 1005.0526
$ws.Range("J58").Value = 2744.5
$ws.Range("K58").Value = 1005.0526
$ws.Range("L58").Value = 2744.5
$ws.Range("M58").Value = -802.0526
$ws.Range("N58").Value = -3150.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3636.3635
$ws.Range("I62").Value = 4400
$ws.Range("J62").Value = 3200
$ws.Range("K62").Value = 4400
$ws.Range("L62").Value = 3200
$ws.Range("M62").Value = -3776
$ws.Range("N62").Value = -4448

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3636.3635
$ws.Range("I65").Value = 4400
$ws.Range("J65").Value = 3200
$ws.Range("K65").Value = 22000
$ws.Range("L65").Value = 16000
$ws.Range("M65").Value = -18880
$ws.Range("N65").Value = -22240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 49346.332
$ws.Range("I132").Value = 2161
$ws.Range("J132").Value = 159445.44
$ws.Range("K132").Value = 6483
$ws.Range("L132").Value = 478336.32
$ws.Range("M132").Value = -3953
$ws.Range("N132").Value = -483396.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 360185.75
$ws.Range("I134").Value = 1033.7354
$ws.Range("J134").Value = 2802419.5
$ws.Range("K134").Value = 3101.2062
$ws.Range("L134").Value = 8407258.5
$ws.Range("M134").Value = -566.2062000000001
$ws.Range("N134").Value = -8412328.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1743
$ws.Range("I136").Value = 1005.0526
$ws.Range("J136").Value = 2744.5
$ws.Range("K136").Value = 3015.1578
$ws.Range("L136").Value = 8233.5
$ws.Range("M136").Value = -465.1578
$ws.Range("N136").Value = -13333.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1418.0869
$ws.Range("J68").Value = 1824.909
$ws.Range("L68").Value = 5474.727000000001
$ws.Range("N68").Value = -7096.727000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1418.0869
$ws.Range("J71").Value = 1824.909
$ws.Range("L71").Value = 16424.181
$ws.Range("N71").Value = -24536.181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2203.0833
$ws.Range("I132").Value = 959.6
$ws.Range("J132").Value = 3091.2856
$ws.Range("K132").Value = 8636.4
$ws.Range("L132").Value = 27821.5704
$ws.Range("M132").Value = -6106.4
$ws.Range("N132").Value = -32881.5704

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 57290.89
$ws.Range("I140").Value = 84184.086
$ws.Range("J140").Value = 3504.5
$ws.Range("K140").Value = 252552.258
$ws.Range("L140").Value = 10513.5
$ws.Range("M140").Value = -247372.258
$ws.Range("N140").Value = -20873.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1285.8823
$ws.Range("I102").Value = 1181.5385
$ws.Range("J102").Value = 1625
$ws.Range("K102").Value = 1181.5385
$ws.Range("L102").Value = 1625
$ws.Range("M102").Value = 440.4614999999999
$ws.Range("N102").Value = -4869

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1650
$ws.Range("I122").Value = 1650
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4950
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2500
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 58827156
$ws.Range("I7").Value = 100002400
$ws.Range("J7").Value = 5376.857
$ws.Range("K7").Value = 100002400
$ws.Range("L7").Value = 5376.857
$ws.Range("M7").Value = -100002288
$ws.Range("N7").Value = -5600.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5287.1333
$ws.Range("I40").Value = 4970.5
$ws.Range("K40").Value = 4970.5
$ws.Range("M40").Value = -4834.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3648.484
$ws.Range("I68").Value = 3559.2273
$ws.Range("J68").Value = 3866.6667
$ws.Range("K68").Value = 3559.2273
$ws.Range("L68").Value = 3866.6667
$ws.Range("M68").Value = -2810.2273
$ws.Range("N68").Value = -5364.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3648.484
$ws.Range("I71").Value = 3559.2273
$ws.Range("J71").Value = 3866.6667
$ws.Range("K71").Value = 17796.1365
$ws.Range("L71").Value = 19333.3335
$ws.Range("M71").Value = -14052.1365
$ws.Range("N71").Value = -26821.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 61044.59
$ws.Range("I122").Value = 68890.53
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 206671.59
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -204221.59
$ws.Range("N122").Value = -11500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 58827156
$ws.Range("I126").Value = 100002400
$ws.Range("J126").Value = 5376.857
$ws.Range("K126").Value = 300007200
$ws.Range("L126").Value = 16130.571
$ws.Range("M126").Value = -300004730
$ws.Range("N126").Value = -21070.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1339.8
$ws.Range("I122").Value = 1339.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4019.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1569.4
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3269728
$ws.Range("I126").Value = 4203222.5
$ws.Range("J126").Value = 2497
$ws.Range("K126").Value = 12609667.5
$ws.Range("L126").Value = 7491
$ws.Range("M126").Value = -12607197.5
$ws.Range("N126").Value = -12431

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3187.1177
$ws.Range("I132").Value = 2683.5715
$ws.Range("J132").Value = 3539.6
$ws.Range("K132").Value = 8050.7145
$ws.Range("L132").Value = 10618.8
$ws.Range("M132").Value = -5520.7145
$ws.Range("N132").Value = -15678.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 257465.84
$ws.Range("I136").Value = 334165.5
$ws.Range("J136").Value = 1800.3334
$ws.Range("K136").Value = 1002496.5
$ws.Range("L136").Value = 5401.0002
$ws.Range("M136").Value = -999946.5
$ws.Range("N136").Value = -10501.0002
